# Apply the "Direct Link Loss Measurement" slide edit:
#  - Resize/reposition the title textbox (Title 2) on slide 10 so it spans
#    the full slide width (Left 0, Width 714pt == 9067800 EMU).
#  - Center-align both paragraphs in that title textbox.
#  - Fix the title text: "Link Loss Direct Loss Measurement (P2P Circuits)"
#    -> "Direct Link Loss Measurement (P2P Circuits)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(1)

# Reposition / resize the title shape.
$sh.Left = 0
$sh.Width = 714

$tr = $sh.TextFrame.TextRange

# Fix the headline text while keeping it inside a single run.
$headline = $tr.Characters(1, 48)
$headline.Text = "Direct Link Loss Measurement (P2P Circuits)"

# Center both paragraphs of the title.
$tr.ParagraphFormat.Alignment = 2
